# The "Do nothing" choice paragraph originally read "--D" + "o nothing"
# split across two runs (with the _GoBack bookmark sandwiched in between,
# left over from whoever last edited the doc in Word). The author forgot to
# update this choice's wording when the "kill Bentley" choices were edited,
# so here we:
#   1. Change the visible text to "--Choose not to kill Bentley".
#   2. Remove the now-redundant trailing "o nothing" run.
#   3. Remove the stray blank paragraph that used to sit right after it.
# The _GoBack bookmark itself is left in place untouched.

$d = $word.ActiveDocument

# Locate the choice paragraph by its current text rather than a hard-coded
# paragraph index, so this keeps working even if the document shifts.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text
    if ($ptxt.Contains("Do nothing")) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not locate the 'Do nothing' choice paragraph"
}

# 1) Grow the "--D" run's text into the full new wording. Setting .Text on
#    the found (small) range only rewrites that run, so the bookmark that
#    follows it in the paragraph is left completely alone.
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$fr = $d.Range($r.Start, $r.End)
$fr.Find.Execute("--D", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fr.Text = "--Choose not to kill Bentley"

# 2) Delete the leftover "o nothing" run (the tail of the original text,
#    which now reads as a dangling continuation after the bookmark).
$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$fr2 = $d.Range($r.Start, $r.End)
$found = $fr2.Find.Execute("o nothing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $fr2.Delete()
}

# 3) Delete the blank paragraph that immediately follows the choice
#    paragraph (Range.Text for an empty paragraph is just the paragraph
#    mark, so check for that rather than an exact empty-string match).
$blankIndex = $targetIndex + 1
if ($blankIndex -le $d.Paragraphs.Count) {
    $blank = $d.Paragraphs.Item($blankIndex)
    if ($blank.Range.Text.Trim().Length -eq 0) {
        $blank.Range.Delete()
    }
}
